# The "ABC分析_客構成" sheet (2nd worksheet) had column H (客単価 / per-
# customer average) computed incorrectly for the top block of rows: it had
# been left holding an older, un-divided figure instead of
# G (平均支払額) divided by C (count_客構成).
# Recompute H = G / C for that block (rows 2-46), skipping any row where
# C is 0 (would divide by zero / already holds a non-numeric placeholder).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

for ($r = 2; $r -le 46; $r++) {
    $count = $ws.Cells.Item($r, 3).Value2
    $avgPay = $ws.Cells.Item($r, 7).Value2
    if ($count -ne 0) {
        $ws.Cells.Item($r, 8).Value = $avgPay / $count
    }
}
